# ---------------------------------------------------------------------------
# Automated BRVM data refresh (GitHub Actions).
#
# Sheet "Recommandations": the per-sector/per-stock stats (days up/down,
# total/last variation, recommendation + strategy tags) are recomputed from
# the latest feed. Rows are re-sorted by the new "Variation Totale (%)"
# ranking, several rows get fresh figures, and the bottom 5 rows (which no
# longer make the cut) are dropped, shrinking the table from A1:G41 to
# A1:G36.
#
# Sheet "Top_YTD": same sector list, refreshed "Progression YTD (%)" values.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 2: BRVM - SERVICES PUBLICS
$ws1.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(2, 2).Value = 0
$ws1.Cells.Item(2, 3).Value = 6
$ws1.Cells.Item(2, 4).Value = 2462.97
$ws1.Cells.Item(2, 5).Value = 103.97
$ws1.Cells.Item(2, 6).Value = "🟡 Observer"
$ws1.Cells.Item(2, 7).Value = "➖ Neutre"

# Row 3: BRVM - AUTRES SECTEURS
$ws1.Cells.Item(3, 1).Value = "BRVM - AUTRES SECTEURS"
$ws1.Cells.Item(3, 2).Value = 0
$ws1.Cells.Item(3, 3).Value = 3
$ws1.Cells.Item(3, 4).Value = 1981.6
$ws1.Cells.Item(3, 5).Value = 701.11
$ws1.Cells.Item(3, 6).Value = "🟡 Observer"
$ws1.Cells.Item(3, 7).Value = "➖ Neutre"

# Row 4: NEI-CEDA CI
$ws1.Cells.Item(4, 1).Value = "NEI-CEDA CI"
$ws1.Cells.Item(4, 2).Value = 0
$ws1.Cells.Item(4, 3).Value = 2
$ws1.Cells.Item(4, 4).Value = 1935
$ws1.Cells.Item(4, 5).Value = 990
$ws1.Cells.Item(4, 6).Value = "🟡 Observer"
$ws1.Cells.Item(4, 7).Value = "➖ Neutre"

# Row 5: BRVM - DISTRIBUTION
$ws1.Cells.Item(5, 1).Value = "BRVM - DISTRIBUTION"
$ws1.Cells.Item(5, 2).Value = 0
$ws1.Cells.Item(5, 3).Value = 3
$ws1.Cells.Item(5, 4).Value = 1474.96
$ws1.Cells.Item(5, 5).Value = 495.75
$ws1.Cells.Item(5, 6).Value = "🟡 Observer"
$ws1.Cells.Item(5, 7).Value = "➖ Neutre"

# Row 6: BRVM - TRANSPORT
$ws1.Cells.Item(6, 1).Value = "BRVM - TRANSPORT"
$ws1.Cells.Item(6, 2).Value = 0
$ws1.Cells.Item(6, 3).Value = 3
$ws1.Cells.Item(6, 4).Value = 1063.54
$ws1.Cells.Item(6, 5).Value = 357.37
$ws1.Cells.Item(6, 6).Value = "🟡 Observer"
$ws1.Cells.Item(6, 7).Value = "➖ Neutre"

# Row 7: BRVM - AGRICULTURE
$ws1.Cells.Item(7, 1).Value = "BRVM - AGRICULTURE"
$ws1.Cells.Item(7, 2).Value = 0
$ws1.Cells.Item(7, 3).Value = 3
$ws1.Cells.Item(7, 4).Value = 986.5
$ws1.Cells.Item(7, 5).Value = 331.17
$ws1.Cells.Item(7, 6).Value = "🟡 Observer"
$ws1.Cells.Item(7, 7).Value = "➖ Neutre"

# Row 8: BRVM - INDUSTRIE  (**)
$ws1.Cells.Item(8, 1).Value = "BRVM - INDUSTRIE  (**)"
$ws1.Cells.Item(8, 2).Value = 0
$ws1.Cells.Item(8, 3).Value = 3
$ws1.Cells.Item(8, 4).Value = 769.53
$ws1.Cells.Item(8, 5).Value = 257.29
$ws1.Cells.Item(8, 6).Value = "🟡 Observer"
$ws1.Cells.Item(8, 7).Value = "➖ Neutre"

# Row 9: BRVM-PRINCIPAL  (**)
$ws1.Cells.Item(9, 1).Value = "BRVM-PRINCIPAL  (**)"
$ws1.Cells.Item(9, 2).Value = 0
$ws1.Cells.Item(9, 3).Value = 3
$ws1.Cells.Item(9, 4).Value = 647.49
$ws1.Cells.Item(9, 5).Value = 216.3
$ws1.Cells.Item(9, 6).Value = "🟡 Observer"
$ws1.Cells.Item(9, 7).Value = "➖ Neutre"

# Row 10: BRVM - CONSOMMATION DE BASE  (**)
$ws1.Cells.Item(10, 1).Value = "BRVM - CONSOMMATION DE BASE  (**)"
$ws1.Cells.Item(10, 2).Value = 0
$ws1.Cells.Item(10, 3).Value = 3
$ws1.Cells.Item(10, 4).Value = 639.88
$ws1.Cells.Item(10, 5).Value = 214.13
$ws1.Cells.Item(10, 6).Value = "🟡 Observer"
$ws1.Cells.Item(10, 7).Value = "➖ Neutre"

# Row 11: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Cells.Item(11, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(11, 2).Value = 0
$ws1.Cells.Item(11, 3).Value = 3
$ws1.Cells.Item(11, 4).Value = 516.16
$ws1.Cells.Item(11, 5).Value = 175.96
$ws1.Cells.Item(11, 6).Value = "🟡 Observer"
$ws1.Cells.Item(11, 7).Value = "➖ Neutre"

# Row 12: BRVM - FINANCES
$ws1.Cells.Item(12, 1).Value = "BRVM - FINANCES"
$ws1.Cells.Item(12, 2).Value = 0
$ws1.Cells.Item(12, 3).Value = 3
$ws1.Cells.Item(12, 4).Value = 436.11
$ws1.Cells.Item(12, 5).Value = 145.91
$ws1.Cells.Item(12, 6).Value = "🟡 Observer"
$ws1.Cells.Item(12, 7).Value = "➖ Neutre"

# Row 13: BRVM - SERVICES FINANCIERS
$ws1.Cells.Item(13, 1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(13, 2).Value = 0
$ws1.Cells.Item(13, 3).Value = 3
$ws1.Cells.Item(13, 4).Value = 428.61
$ws1.Cells.Item(13, 5).Value = 143.4
$ws1.Cells.Item(13, 6).Value = "🟡 Observer"
$ws1.Cells.Item(13, 7).Value = "➖ Neutre"

# Row 14: BRVM-PRESTIGE
$ws1.Cells.Item(14, 1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(14, 2).Value = 0
$ws1.Cells.Item(14, 3).Value = 3
$ws1.Cells.Item(14, 4).Value = 420.33
$ws1.Cells.Item(14, 5).Value = 140.98
$ws1.Cells.Item(14, 6).Value = "🟡 Observer"
$ws1.Cells.Item(14, 7).Value = "➖ Neutre"

# Row 15: BRVM - INDUSTRIELS
$ws1.Cells.Item(15, 1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(15, 2).Value = 0
$ws1.Cells.Item(15, 3).Value = 3
$ws1.Cells.Item(15, 4).Value = 394.76
$ws1.Cells.Item(15, 5).Value = 132.83
$ws1.Cells.Item(15, 6).Value = "🟡 Observer"
$ws1.Cells.Item(15, 7).Value = "➖ Neutre"

# Row 16: BRVM - ENERGIE
$ws1.Cells.Item(16, 1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(16, 2).Value = 0
$ws1.Cells.Item(16, 3).Value = 3
$ws1.Cells.Item(16, 4).Value = 328.57
$ws1.Cells.Item(16, 5).Value = 109.84
$ws1.Cells.Item(16, 6).Value = "🟡 Observer"
$ws1.Cells.Item(16, 7).Value = "➖ Neutre"

# Row 17: BRVM - TELECOMMUNICATIONS
$ws1.Cells.Item(17, 1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(17, 2).Value = 0
$ws1.Cells.Item(17, 3).Value = 3
$ws1.Cells.Item(17, 4).Value = 280.34
$ws1.Cells.Item(17, 5).Value = 93.28
$ws1.Cells.Item(17, 6).Value = "🟡 Observer"
$ws1.Cells.Item(17, 7).Value = "➖ Neutre"

# Row 18: SETAO CI (STAC)
$ws1.Cells.Item(18, 1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(18, 2).Value = 3
$ws1.Cells.Item(18, 3).Value = 0
$ws1.Cells.Item(18, 4).Value = 21.76
$ws1.Cells.Item(18, 5).Value = 7.14
$ws1.Cells.Item(18, 6).Value = "🟢 Achat"
$ws1.Cells.Item(18, 7).Value = "✅ Renforcer"

# Row 19: LOTERIE NATIONALE DU BENIN (LNBB)
$ws1.Cells.Item(19, 1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Cells.Item(19, 2).Value = 2
$ws1.Cells.Item(19, 3).Value = 0
$ws1.Cells.Item(19, 4).Value = 11.68
$ws1.Cells.Item(19, 5).Value = 6.55
$ws1.Cells.Item(19, 6).Value = "🟡 Observer"
$ws1.Cells.Item(19, 7).Value = "➖ Neutre"

# Row 20: NEI-CEDA CI (NEIC)
$ws1.Cells.Item(20, 1).Value = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(20, 2).Value = 2
$ws1.Cells.Item(20, 3).Value = 0
$ws1.Cells.Item(20, 4).Value = 10.88
$ws1.Cells.Item(20, 5).Value = 6.12
$ws1.Cells.Item(20, 6).Value = "🟡 Observer"
$ws1.Cells.Item(20, 7).Value = "➖ Neutre"

# Row 21: SICABLE CI (CABC)
$ws1.Cells.Item(21, 1).Value = "SICABLE CI (CABC)"
$ws1.Cells.Item(21, 2).Value = 1
$ws1.Cells.Item(21, 3).Value = 0
$ws1.Cells.Item(21, 4).Value = 7.5
$ws1.Cells.Item(21, 5).Value = 7.5
$ws1.Cells.Item(21, 6).Value = "🟡 Observer"
$ws1.Cells.Item(21, 7).Value = "➖ Neutre"

# Row 22: ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)
$ws1.Cells.Item(22, 1).Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Cells.Item(22, 2).Value = 1
$ws1.Cells.Item(22, 3).Value = 1
$ws1.Cells.Item(22, 4).Value = 4.15
$ws1.Cells.Item(22, 5).Value = -3.1
$ws1.Cells.Item(22, 6).Value = "🟡 Observer"
$ws1.Cells.Item(22, 7).Value = "👀 À surveiller"

# Row 23: BANK OF AFRICA NG (BOAN)
$ws1.Cells.Item(23, 1).Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Cells.Item(23, 2).Value = 2
$ws1.Cells.Item(23, 3).Value = 1
$ws1.Cells.Item(23, 4).Value = 2.45
$ws1.Cells.Item(23, 5).Value = 4.13
$ws1.Cells.Item(23, 6).Value = "🟡 Observer"
$ws1.Cells.Item(23, 7).Value = "👀 À surveiller"

# Row 24: BANK OF AFRICA BF (BOABF)
$ws1.Cells.Item(24, 1).Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Cells.Item(24, 2).Value = 1
$ws1.Cells.Item(24, 3).Value = 1
$ws1.Cells.Item(24, 4).Value = 1.31
$ws1.Cells.Item(24, 5).Value = 4.86
$ws1.Cells.Item(24, 6).Value = "🟡 Observer"
$ws1.Cells.Item(24, 7).Value = "👀 À surveiller"

# Row 25: TOTAL
$ws1.Cells.Item(25, 1).Value = "TOTAL"
$ws1.Cells.Item(25, 2).Value = 0
$ws1.Cells.Item(25, 3).Value = 3
$ws1.Cells.Item(25, 4).Value = 0
$ws1.Cells.Item(25, 5).Value = 0
$ws1.Cells.Item(25, 6).Value = "🟡 Observer"
$ws1.Cells.Item(25, 7).Value = "➖ Neutre"

# Row 26: CFAO MOTORS CI (CFAC)
$ws1.Cells.Item(26, 1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(26, 2).Value = 1
$ws1.Cells.Item(26, 3).Value = 1
$ws1.Cells.Item(26, 4).Value = -0.24
$ws1.Cells.Item(26, 5).Value = 3.21
$ws1.Cells.Item(26, 6).Value = "🟡 Observer"
$ws1.Cells.Item(26, 7).Value = "👀 À surveiller"

# Row 27: SAPH CI (SPHC)
$ws1.Cells.Item(27, 1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(27, 2).Value = 1
$ws1.Cells.Item(27, 3).Value = 1
$ws1.Cells.Item(27, 4).Value = -0.73
$ws1.Cells.Item(27, 5).Value = 4.93
$ws1.Cells.Item(27, 6).Value = "🟡 Observer"
$ws1.Cells.Item(27, 7).Value = "👀 À surveiller"

# Row 28: NESTLE CI (NTLC)
$ws1.Cells.Item(28, 1).Value = "NESTLE CI (NTLC)"
$ws1.Cells.Item(28, 2).Value = 0
$ws1.Cells.Item(28, 3).Value = 1
$ws1.Cells.Item(28, 4).Value = -0.96
$ws1.Cells.Item(28, 5).Value = -0.96
$ws1.Cells.Item(28, 6).Value = "🟡 Observer"
$ws1.Cells.Item(28, 7).Value = "➖ Neutre"

# Row 29: SUCRIVOIRE (SCRC)
$ws1.Cells.Item(29, 1).Value = "SUCRIVOIRE (SCRC)"
$ws1.Cells.Item(29, 2).Value = 1
$ws1.Cells.Item(29, 3).Value = 2
$ws1.Cells.Item(29, 4).Value = -1.26
$ws1.Cells.Item(29, 5).Value = -0.9
$ws1.Cells.Item(29, 6).Value = "🟡 Observer"
$ws1.Cells.Item(29, 7).Value = "👀 À surveiller"

# Row 30: SOCIETE IVOIRIENNE DE BANQUE  (SIBC)
$ws1.Cells.Item(30, 1).Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Cells.Item(30, 2).Value = 0
$ws1.Cells.Item(30, 3).Value = 1
$ws1.Cells.Item(30, 4).Value = -1.5
$ws1.Cells.Item(30, 5).Value = -1.5
$ws1.Cells.Item(30, 6).Value = "🟡 Observer"
$ws1.Cells.Item(30, 7).Value = "➖ Neutre"

# Row 31: BANK OF AFRICA ML (BOAM)
$ws1.Cells.Item(31, 1).Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Cells.Item(31, 2).Value = 0
$ws1.Cells.Item(31, 3).Value = 1
$ws1.Cells.Item(31, 4).Value = -2.68
$ws1.Cells.Item(31, 5).Value = -2.68
$ws1.Cells.Item(31, 6).Value = "🟡 Observer"
$ws1.Cells.Item(31, 7).Value = "➖ Neutre"

# Row 32: SOLIBRA CI (SLBC)
$ws1.Cells.Item(32, 1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(32, 2).Value = 0
$ws1.Cells.Item(32, 3).Value = 1
$ws1.Cells.Item(32, 4).Value = -2.74
$ws1.Cells.Item(32, 5).Value = -2.74
$ws1.Cells.Item(32, 6).Value = "🟡 Observer"
$ws1.Cells.Item(32, 7).Value = "➖ Neutre"

# Row 33: NSIA BANQUE COTE D'IVOIRE (NSBC)
$ws1.Cells.Item(33, 1).Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Cells.Item(33, 2).Value = 0
$ws1.Cells.Item(33, 3).Value = 1
$ws1.Cells.Item(33, 4).Value = -3.51
$ws1.Cells.Item(33, 5).Value = -3.51
$ws1.Cells.Item(33, 6).Value = "🟡 Observer"
$ws1.Cells.Item(33, 7).Value = "➖ Neutre"

# Row 34: CIE CI (CIEC)
$ws1.Cells.Item(34, 1).Value = "CIE CI (CIEC)"
$ws1.Cells.Item(34, 2).Value = 0
$ws1.Cells.Item(34, 3).Value = 1
$ws1.Cells.Item(34, 4).Value = -4.09
$ws1.Cells.Item(34, 5).Value = -4.09
$ws1.Cells.Item(34, 6).Value = "🟡 Observer"
$ws1.Cells.Item(34, 7).Value = "➖ Neutre"

# Row 35: ECOBANK TRANS. INCORP. TG (ETIT)
$ws1.Cells.Item(35, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(35, 2).Value = 0
$ws1.Cells.Item(35, 3).Value = 1
$ws1.Cells.Item(35, 4).Value = -4.35
$ws1.Cells.Item(35, 5).Value = -4.35
$ws1.Cells.Item(35, 6).Value = "🟡 Observer"
$ws1.Cells.Item(35, 7).Value = "➖ Neutre"

# Row 36: VIVO ENERGY CI (SHEC)
$ws1.Cells.Item(36, 1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(36, 2).Value = 0
$ws1.Cells.Item(36, 3).Value = 1
$ws1.Cells.Item(36, 4).Value = -6.74
$ws1.Cells.Item(36, 5).Value = -6.74
$ws1.Cells.Item(36, 6).Value = "🟡 Observer"
$ws1.Cells.Item(36, 7).Value = "➖ Neutre"

# Drop the bottom 5 rows that no longer make the cut -> table becomes A1:G36.
$ws1.Range("A37:G41").Delete()

# Refresh "Top_YTD" progression figures (sector/title order is unchanged).
$ws2 = $wb.Worksheets.Item("Top_YTD")
$ws2.Cells.Item(2, 2).Value = 464721.14
$ws2.Cells.Item(3, 2).Value = 43774.32
$ws2.Cells.Item(4, 2).Value = 20609.64
$ws2.Cells.Item(5, 2).Value = 11290.5
$ws2.Cells.Item(6, 2).Value = 9289.01
$ws2.Cells.Item(7, 2).Value = 7785.7
$ws2.Cells.Item(8, 2).Value = 4431.17
$ws2.Cells.Item(9, 2).Value = 3050.35
$ws2.Cells.Item(10, 2).Value = 2975.04
$ws2.Cells.Item(11, 2).Value = 1913.1
